$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the oldest data row (old row 2), shifting all subsequent rows up by one.
$ws.Rows(2).Delete()

# Update the recalculated year-over-year forecast values (column E) for all
# remaining data rows (now rows 2-18) to reflect the corrected forecaster output.
$ws.Range("E2").Value  = -4.700509864312973
$ws.Range("E3").Value  = -0.01655958389530365
$ws.Range("E4").Value  = 3.579142225970444
$ws.Range("E5").Value  = -0.289184878867832
$ws.Range("E6").Value  = 5.963492031746176
$ws.Range("E7").Value  = 7.523777575896196
$ws.Range("E8").Value  = 2.532215190177589
$ws.Range("E9").Value  = 2.051185924063259
$ws.Range("E10").Value = 0.4575538530338541
$ws.Range("E11").Value = 2.600569166164624
$ws.Range("E12").Value = 3.605726003451304
$ws.Range("E13").Value = 3.490656491795074
$ws.Range("E14").Value = -2.347097924577757
$ws.Range("E15").Value = -0.1803381976702711
$ws.Range("E16").Value = -1.152671696465724
$ws.Range("E17").Value = -2.785556326028149
$ws.Range("E18").Value = -2.452009576682213
